# ValueSet-VSMotivoSolicitudHosp.xlsx — "version final sin errores"
#
# Two changes on the "Metadata" sheet:
#   1) Bump the Version value (B3) from "0.4.0" to "0.7.0".
#   2) Remove the "Jurisdiction" / "Chile" row entirely (row 11), which
#      shifts every row below it up by one (Description, Purpose,
#      Copyright, Immutable all move up one row) and shrinks the used
#      range from A1:B15 to A1:B14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1) Version 0.4.0 -> 0.7.0
$ws.Range("B3").Value = "0.7.0"

# 2) Delete the whole "Jurisdiction" / "Chile" row
$ws.Rows.Item(11).Delete()
